$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 7 (old Layout1 thumbnail, Layout2 view/thumbnail rows)
$ws.Rows("4:7").Delete()

# Clear the old wrap-text formatting on E2 before rewriting the row
$ws.Range("E2").Clear()
$ws.Rows("2:2").AutoFit()

# Rewrite row 2: scaffold_context_info.json entry
$ws.Range("A2").Value = "scaffold_context_info.json"
$ws.Range("B2").Value = "application/x.vnd.abi.context-information+json"
$ws.Range("C2").Value = '{"version": "0.2.0", "id": "sparc.science.annotation_terms"}'
$ws.Range("D2").Value = "rat_brainstem_metadata.json"

# Rewrite row 3: rat_brainstem_metadata.json entry
$ws.Range("A3").Value = "rat_brainstem_metadata.json"
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "scaffold_context_info.json"

# Adjust column widths (A and D shrink to match the shorter filenames,
# E shrinks since it now only holds a short cross-reference)
$ws.Columns("A").ColumnWidth = 23
$ws.Columns("D").ColumnWidth = 23
$ws.Columns("E").ColumnWidth = 8.5

# Select E3 as active cell
$ws.Range("E3").Select()
